# EPBDS-2195 Rem('%') operator for BigDecimalValue, BigIntegerValue
# Add two new test blocks (BigIntegerValue, BigDecimalValue) to the Rem operator
# test sheet, mirroring the existing DoubleValue block (E56:F62).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Block 1: testRemBigIntegerValue, rows 66-72 -----------------------
# Merge the header cells *before* copying formats over them: merging an
# already-formatted range re-derives border styles, so doing it first
# (while the cells are still blank/default) and then pasting the real
# formatting on top keeps the style indices identical to the source block.
$ws.Range("E66:F66").Merge()

# Copy formatting (fills/fonts/borders/row heights) from the DoubleValue
# block so the new block looks identical to the existing ones.
$srcBlock = $ws.Range("E56:F62")
$srcBlock.Copy()
$ws.Range("E66:F72").PasteSpecial(-4122)

$ws.Range("E66").Value = "Rules String testRemBigIntegerValue(BigIntegerValue v1, BigIntegerValue v2)"

$ws.Range("E67").Value = "C1"
$ws.Range("F67").Value = "RET1"

$ws.Range("E68").Value = "v1 % v2 == 0"

$ws.Range("E69").Value = "boolean"

$ws.Range("E70").Value = "Value1 - Value2"
$ws.Range("F70").Value = "Result"

$ws.Range("E71").Value = "Yes"
$ws.Range("F71").Value = "passed"

$ws.Range("E72").Value = "No"
$ws.Range("F72").Value = "not passed"

$ws.Rows.Item(66).RowHeight = 17.25
$ws.Rows.Item(70).RowHeight = 17.25
$ws.Rows.Item(71).RowHeight = 17.25
$ws.Rows.Item(72).RowHeight = 17.25

# --- Block 2: testRemBigDecimalValue, rows 76-82 -----------------------
$ws.Range("E76:F76").Merge()

$srcBlock2 = $ws.Range("E56:F62")
$srcBlock2.Copy()
$ws.Range("E76:F82").PasteSpecial(-4122)

$ws.Range("E76").Value = "Rules String testRemBigDecimalValue(BigDecimalValue v1, BigDecimalValue v2)"

$ws.Range("E77").Value = "C1"
$ws.Range("F77").Value = "RET1"

$ws.Range("E78").Value = "v1 % v2 == 0"

$ws.Range("E79").Value = "boolean"

$ws.Range("E80").Value = "Value1 - Value2"
$ws.Range("F80").Value = "Result"

$ws.Range("E81").Value = "Yes"
$ws.Range("F81").Value = "passed"

$ws.Range("E82").Value = "No"
$ws.Range("F82").Value = "not passed"

$ws.Rows.Item(76).RowHeight = 17.25
$ws.Rows.Item(80).RowHeight = 17.25
$ws.Rows.Item(81).RowHeight = 17.25
$ws.Rows.Item(82).RowHeight = 17.25

# --- View state: scroll / selection -------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E78").Select()
